# Applies the cryptos-list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.097.97"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").Value = "1.850.14"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'0.6928"
$ws.Range("E5").Value = "  -5.40%  "
$ws.Range("D6").Value = "'237.89"
$ws.Range("E6").Value = "  -1.39%  "
$ws.Range("D7").Value = "'0.9998"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.07743"
$ws.Range("E8").Value = "  +8.54%  "
$ws.Range("D9").Value = "'0.3037"
$ws.Range("E9").Value = "  -3.14%  "
$ws.Range("D10").Value = "'23.27"
$ws.Range("E10").Value = "  -4.65%  "
$ws.Range("D11").Value = "'0.08113"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").Value = "1.851.77"
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").Value = "'5.205"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "'89.07"
$ws.Range("D16").Value = "29.102.38"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").Value = "'5.738"
$ws.Range("E17").Value = "  -4.51%  "
$ws.Range("D18").Value = "'0.000007787"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").Value = "'13.19"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("D20").Value = "'235.50"
$ws.Range("E20").Value = "  -5.11%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").Value = "2.099.16"
$ws.Range("E22").Value = "  -2.23%  "
$ws.Range("D23").Value = "'0.9998"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "'7.600"
$ws.Range("E24").Value = "  -2.07%  "
$ws.Range("D25").Value = "'8.975"
$ws.Range("E25").Value = "  -2.47%  "
$ws.Range("D26").Value = "'160.90"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("D27").Value = "'0.1429"
$ws.Range("E27").Value = "  -7.38%  "
$ws.Range("D28").Value = "'18.05"
$ws.Range("E28").Value = "  -2.65%  "
$ws.Range("D29").Value = "'1.975"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("D30").Value = "'1.400"
$ws.Range("E30").Value = "  -3.16%  "
$ws.Range("D31").Value = "'4.489"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").Value = "'4.016"
$ws.Range("E33").Value = "  -3.97%  "
$ws.Range("D34").Value = "'0.05230"
$ws.Range("E34").Value = "  -1.33%  "
$ws.Range("D35").Value = "'1.179"
$ws.Range("E35").Value = "  -4.21%  "
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "'1.026"
$ws.Range("E36").Value = "  +2.49%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.7019"
$ws.Range("E37").Value = "  -5.22%  "
$ws.Range("D38").Value = "'2.651"
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("E39").Value = "  -4.19%  "
$ws.Range("D40").Value = "'2.678"
$ws.Range("E40").Value = "  -2.02%  "
$ws.Range("D41").Value = "'0.9158"
$ws.Range("E41").Value = "  +5.73%  "
$ws.Range("D42").Value = "1.089.17"
$ws.Range("E42").Value = "  +4.53%  "
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").Value = "'0.4260"
$ws.Range("E44").Value = "  -4.40%  "
$ws.Range("D45").Value = "'70.65"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").Value = "'0.9996"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "'103.10"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("D49").Value = "1.997.02"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("D50").Value = "'9.149"
$ws.Range("E50").Value = "  -3.63%  "
$ws.Range("D51").Value = "'6.973"
$ws.Range("E51").Value = "  -6.10%  "
